# Weekly refresh of the "Arándano (blue)" price series at Vega Modelo de
# Temuco: a new day's worth of records is inserted at the top of the data
# block (row 4) and every existing record shifts down by one row, with the
# last existing row re-appearing one row further down (old row 101 -> new
# row 102).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually vary row-to-row in this data block; A,B,C,E,F,G,H,I,J,K
# are constant across the whole sheet and don't need touching.
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$colIdx = @{ "D"=4; "L"=12; "M"=13; "N"=14; "O"=15; "P"=16; "Q"=17; "R"=18; "S"=19; "T"=20 }

$firstDataRow = 4
$lastDataRow = 101

# 1) Snapshot the existing data rows (4..101) before overwriting anything,
#    since the shift-down would otherwise clobber values we still need to
#    read.
$savedRows = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $colIdx[$c]).Value2
    }
    $savedRows[$r] = $rowData
}

# 2) Shift every saved row down by one: new row (r+1) gets old row r's data.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $destRow = $r + 1
    $rowData = $savedRows[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $colIdx[$c]).Value2 = $rowData[$c]
    }
}

# 3) Populate the brand-new record at row 4 with this week's reading.
$ws.Cells.Item($firstDataRow, $colIdx["D"]).Value2 = 44882
$ws.Cells.Item($firstDataRow, $colIdx["L"]).Value2 = "Primera"
$ws.Cells.Item($firstDataRow, $colIdx["M"]).Value2 = 300
$ws.Cells.Item($firstDataRow, $colIdx["N"]).Value2 = 3200
$ws.Cells.Item($firstDataRow, $colIdx["O"]).Value2 = 3200
$ws.Cells.Item($firstDataRow, $colIdx["P"]).Value2 = 3200
$ws.Cells.Item($firstDataRow, $colIdx["Q"]).Value2 = "`$/kilo"
$ws.Cells.Item($firstDataRow, $colIdx["R"]).Value2 = "Región del Maule"
$ws.Cells.Item($firstDataRow, $colIdx["S"]).Value2 = 3200
$ws.Cells.Item($firstDataRow, $colIdx["T"]).Value2 = 1

# 4) Fill in the constant columns (A,B,C,E..K) for the newly created row 102,
#    copied from the row directly above it (these are identical for every
#    record in this block).
$constCols = @{ "A"=1; "B"=2; "C"=3; "E"=5; "F"=6; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11 }
foreach ($cc in $constCols.Keys) {
    $ci = $constCols[$cc]
    $ws.Cells.Item(102, $ci).Value2 = $ws.Cells.Item(101, $ci).Value2
}

# Keep the sheet's declared dimension honest.
$ws.Range("D102").NumberFormat = $ws.Range("D101").NumberFormat
